# Move the new analysis sheet after Sheet1, populate it, and update
# Sheet1's view state (tab/selection) to reflect the new active sheet.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# --- create the new worksheet right after Sheet1 ---------------------------
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)
$ws.Name = "AXI2MEM_analysis"

# --- tile size header --------------------------------------------------
$ws.Range("B2").Value = "tile size:"
$ws.Range("B2").HorizontalAlignment = -4152
$ws.Range("C2").Value = 4
$ws.Range("C2").HorizontalAlignment = -4131

# --- pxl / tiles table ---------------------------------------------------
$ws.Range("D5").Value = "pxl"
$ws.Range("E5").Value = "tiles"
$ws.Range("C5:E5").HorizontalAlignment = -4108
$ws.Range("C5:E5").Font.Bold = $true
$ws.Range("C5:E5").Borders.LineStyle = 1

$ws.Range("C6").Value = "height"
$ws.Range("C7").Value = "width"
$ws.Range("C8").Value = "total"
$ws.Range("C6:C8").HorizontalAlignment = -4152
$ws.Range("C6:C8").Font.Bold = $true
$ws.Range("C6:C8").Borders.LineStyle = 1

$ws.Range("D6").Value = 480
$ws.Range("D7").Value = 640
$ws.Range("D8").Formula = "=D7*D6"
$ws.Range("E6").Formula = "=D6/C2"
$ws.Range("E7").Formula = "=D7/C2"
$ws.Range("E8").Formula = "=E7*E6"
$ws.Range("D6:E8").HorizontalAlignment = -4108
$ws.Range("D6:E8").Borders.LineStyle = 1

# --- explanatory text ------------------------------------------------------
$ws.Range("C11").Value = "Each AXI address is " + [char]0x00BD + " of a frame buffer line"
$ws.Range("C12").Value = "To convert address, divide it by 4."
$ws.Range("C13").Value = " If address is odd, it is upper half (bit 0 -> 23 of AXI map to bit 24 " + [char]0x2192 + " 47 of frame buff)"
$ws.Range("C14").Value = " If address is even, it is lower half (bit 0 -> 23 of AXI map to bit 0 " + [char]0x2192 + " 23 of frame buff)"
$ws.Range("C13:C14").Font.Bold = $true

# --- AXI / intermediate / MEM table ----------------------------------------
$ws.Range("C17").Value = "AXI"
$ws.Range("D17").Value = "intermediate"
$ws.Range("E17").Value = "MEM"
$ws.Range("C17:E17").HorizontalAlignment = -4108
$ws.Range("C17:E17").Font.Bold = $true
$ws.Range("C17:E17").Borders.LineStyle = 1

$ws.Range("C18").Value = 0
$ws.Range("C19").Value = 4
$ws.Range("C20").Value = 8
$ws.Range("C21").Formula = "=C20+4"
$ws.Range("C22").Formula = "=C21+4"
$ws.Range("C23").Formula = "=C22+4"
$ws.Range("C24").Formula = "=C23+4"
$ws.Range("C25").Formula = "=C24+4"
$ws.Range("C26").Formula = "=C25+4"
$ws.Range("C27").Formula = "=C26+4"
$ws.Range("C28").Formula = "=C27+4"
$ws.Range("C29").Formula = "=C28+4"
$ws.Range("C30").Formula = "=C29+4"

foreach ($r in 18..30) {
    $ws.Range("D$r").Formula = "=C$r/4"
    $ws.Range("E$r").Formula = "=IF((ISODD(D$r)), (D$r-1)/2, D$r/2)"
}

$ws.Range("C31").Value = "…"
$ws.Range("D31").Value = "…"
$ws.Range("E31").Value = "…"

$ws.Range("C32").Formula = "=(E8-1)*2*4"
$ws.Range("D32").Formula = "=C32/4"
$ws.Range("E32").Formula = "=IF((ISODD(D32)), (D32-1)/2, D32/2)"

$ws.Range("C18:E32").HorizontalAlignment = -4108
$ws.Range("C18:E32").Borders.LineStyle = 1

# --- sheet view / selection state -------------------------------------------
$sheet1.Range("M29").Select()
$ws.Range("J20").Select()
$ws.Activate()
